$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 112182656
$ws.Range("B2").Value = 77039
$ws.Range("E2").Value = 6437
$ws.Range("F2").Value = "Blanksvart spiklav"
$ws.Range("G2").Value = "Calicium denigratum"
$ws.Range("H2").Value = "(Vain.) Tibell"
$ws.Range("Q2").Value = 358671
$ws.Range("R2").Value = 6874549

# Row 3 updates
$ws.Range("A3").Value = 112182610
$ws.Range("B3").Value = 77388
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 358653
$ws.Range("R3").Value = 6874558
